# Insert a new data row at row 173 (this pushes the existing rows 173-257
# down to 174-258, matching the target dimension A1:R258) and populate the
# new row with its values. Columns that are identical to the row that used
# to occupy position 173 (A,B,C,E,F,G,H,I,N,O,Q,R) are simply carried over;
# D,J,K,L,M,P get the new figures from this week's price report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(173).Insert()

$ws.Cells.Item(173, 1).Value = 7
$ws.Cells.Item(173, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(173, 3).Value = "Ñuble"
$ws.Cells.Item(173, 4).Value = 44466
$ws.Cells.Item(173, 5).Value = 16
$ws.Cells.Item(173, 6).Value = 100114014
$ws.Cells.Item(173, 7).Value = "Betarraga"
$ws.Cells.Item(173, 8).Value = "Sin especificar"
$ws.Cells.Item(173, 9).Value = "Primera"
$ws.Cells.Item(173, 10).Value = 300
$ws.Cells.Item(173, 11).Value = 750
$ws.Cells.Item(173, 12).Value = 800
$ws.Cells.Item(173, 13).Value = 775
$ws.Cells.Item(173, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(173, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(173, 16).Value = 155
$ws.Cells.Item(173, 17).Value = 5
$ws.Cells.Item(173, 18).Value = "Hortaliza"
